# Weekly driver report update for 2025-04-21
# Update the "Bad Drivers" table (rows 4-7) with refreshed sample data and
# recompute the Totals row (row 9) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.1.3
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.1.3"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 98.3

# Row 5: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 98.40000000000001

# Row 6: Intel(R) Wi-Fi 6 AX201 160MHz - 22.20.0.6
$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.20.0.6"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 98.5

# Row 7: Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6"
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 135
$ws.Range("D7").Value = 98.8

# Row 9: Totals
$ws.Range("B9").Value = 19
$ws.Range("C9").Value = 187
